# "Update Data Sources from LFX"
#
# 1) Re-point every table in the deck from the old local table style
#    ({483223DA-87E1-4EE6-8DA5-15799A8F6AF1} "Table_0") to the new style
#    ({8FA931E9-57A9-4F2B-8CB2-13A407556EF8}).
# 2) Swap the two theme color palettes ("LF Energy Theme 2023" / Geometric
#    and "Simple Light") that are attached to the deck's two slide masters.

$p = $ppt.ActivePresentation

$newTableStyle = "{8FA931E9-57A9-4F2B-8CB2-13A407556EF8}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# Swap the "Geometric" / "LF Energy Theme 2023" color scheme into the
# presentation's (reachable) theme color scheme - RGB values are passed as
# 0xBBGGRR (VBA RGB() long layout) so that they serialize as the intended
# RRGGBB srgbClr values.
$design = $p.Designs.Item(1)
$tcs = $design.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x222222   # dk1       -> 222222
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1       -> FFFFFF
$tcs.Colors(3).RGB  = 0x434343   # dk2       -> 434343
$tcs.Colors(4).RGB  = 0x999999   # lt2       -> 999999
$tcs.Colors(5).RGB  = 0x783700   # accent1   -> 003778
$tcs.Colors(6).RGB  = 0xFF9400   # accent2   -> 0094FF
$tcs.Colors(7).RGB  = 0xE71D5B   # accent3   -> 5B1DE7
$tcs.Colors(8).RGB  = 0xE2E212   # accent4   -> 12E2E2
$tcs.Colors(9).RGB  = 0xAA00FF   # accent5   -> FF00AA
$tcs.Colors(10).RGB = 0x1FDEAC   # accent6   -> ACDE1F
$tcs.Colors(11).RGB = 0xCC7700   # hlink     -> 0077CC
$tcs.Colors(12).RGB = 0x9262F0   # folHlink  -> F06292
